$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# Insert 7 new rows before the current row 4 (old row 4 -> becomes row 11)
$ws.Range("A4:K10").EntireRow.Insert()

# --- Fill column A (test case names) first, rows 5..10 then 4, matching the
#     shared-string insertion order recorded in the target file ---
$ws.Cells.Item(5, 1).Value  = "Login_Error_On_Wrong_Userid"
$ws.Cells.Item(6, 1).Value  = "Login_Error_On_Wrong_Password"
$ws.Cells.Item(7, 1).Value  = "Login_Error_On_Wrong_Token"
$ws.Cells.Item(8, 1).Value  = "Login_Error_On_AppAccess_Restricted"
$ws.Cells.Item(9, 1).Value  = "Login_Error_On_Member_LoginRestricted"
$ws.Cells.Item(10, 1).Value = "Login_Error_On_Member_Deleted"
$ws.Cells.Item(4, 1).Value  = "Login_Error_On_No_Data"

# C6 holds "1234" as text (new shared string)
$ws.Cells.Item(6, 3).Value = "'1234"

# K1 header
$ws.Cells.Item(1, 11).Value = "ExpectedError"

# K7 / D7
$ws.Cells.Item(7, 11).Value = "Token authentication failure."
$ws.Cells.Item(7, 4).Value  = "XT131"

# K5 (new), K6 (reuses K5's new string)
$ws.Cells.Item(5, 11).Value = "ERROR: Incorrect credentials."
$ws.Cells.Item(6, 11).Value = "ERROR: Incorrect credentials."

# B5 test1234 (new)
$ws.Cells.Item(5, 2).Value = "test1234"

# --- Remaining cells reuse existing shared strings, order does not matter ---
$ws.Cells.Item(6, 2).Value = "test123"
$ws.Cells.Item(5, 3).Value = "'123"
$ws.Cells.Item(7, 3).Value = "'123"
$ws.Cells.Item(5, 4).Value = "XT13"
$ws.Cells.Item(6, 4).Value = "XT13"
$ws.Cells.Item(7, 2).Value = "test123"

# C4 is empty but carries the quote-prefix style (s="1") like C2/C3/C5/C6/C7 -
# copy format only (no value) from C2
$ws.Range("C2").Copy()
$ws.Range("C4").PasteSpecial(-4122)

# Rows 8-10 must not inherit the quote-prefix style that Insert() copied down
# from row 3 - remove the stray empty styled cells entirely.
$ws.Range("C8:C10").Clear()

# Column widths (target stored widths are 38.85546875 / 27.28515625; this
# runtime quantizes ColumnWidth to 1/6-character steps, so we pick the input
# that lands on the closest reachable stored width)
$ws.Columns.Item(1).ColumnWidth = 38
$ws.Columns.Item(11).ColumnWidth = 26.5

# Selection
$ws.Range("B5").Select()
